$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values for columns D (open_price), E (close_price), F (high_price),
# G (low_price), H (shares_outstanding), I (fixed_ticker) for rows 2-14.
# Column I is set to "HUT" for every row (previously each row referenced a
# different, now-removed ticker string).
$rows = @(
    @{ Row = 2;  D = 6.800000190734863;  E = 10.60000038146973; F = 11.85000038146973; G = 6.349999904632568;  H = 105527928 },
    @{ Row = 3;  D = 9.199999809265137;  E = 11.19999980926514; F = 12.39999961853027; G = 8.350000381469727;  H = 105527928 },
    @{ Row = 4;  D = 4.264999866485596;  E = 9.649999618530272; F = 10.60000038146973; G = 4.054999828338623;  H = 105527928 },
    @{ Row = 5;  D = 8.949999809265137;  E = 8.949999809265137; F = 11.19999980926514; G = 7.949999809265137;  H = 105527928 },
    @{ Row = 6;  D = 16.35000038146973;  E = 17.79999923706055; F = 22.75;             G = 16.25;              H = 105527928 },
    @{ Row = 7;  D = 10.5;               E = 10.85000038146973; F = 12.32499980926514; G = 8.5;                H = 105527928 },
    @{ Row = 8;  D = 15.0600004196167;   E = 7.46999979019165;  F = 15.06999969482422; G = 6.179999828338623;  H = 105527928 },
    @{ Row = 9;  D = 11.02999973297119;  E = 7.860000133514404; F = 11.6899995803833;  G = 6.949999809265137;  H = 105527928 },
    @{ Row = 10; D = 15.30000019073486;  E = 14.64000034332275; F = 21.0979995727539;  G = 13.85999965667725;  H = 105527928 },
    @{ Row = 11; D = 12.27000045776367;  E = 15.78999996185303; F = 18.38500022888184; G = 10.70199966430664;  H = 105527928 },
    @{ Row = 12; D = 21.10000038146973;  E = 21.69000053405762; F = 29.28000068664551; G = 18.64500045776367;  H = 105527928 },
    @{ Row = 13; D = 11.86999988555908;  E = 12.3100004196167;  F = 13.80000019073486; G = 10.03999996185303;  H = 105527928 },
    @{ Row = 14; D = 18.07999992370605;  E = 21.22999954223633; F = 23.55999946594238; G = 18.03000068664551;  H = 105527928 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
    $ws.Cells.Item($r, 9).Value = "HUT"
}
